$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tb = $s.Shapes.Item(3)

# The textbox currently holds the text "An image" split across three runs
# ("An", " ", "image"). Re-assigning the same concatenated text is a no-op
# for the engine's diffing, so first set it to a distinct placeholder value
# to force the run structure to be rebuilt, then set the final text. This
# collapses the three runs into a single run with no explicit formatting.
$tb.TextFrame.TextRange.Text = "__tmp__"
$tb.TextFrame.TextRange.Text = "An image"
